$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Realtek RTL8811AU ... (A3/B3 unchanged) ---
$ws.Range("C3").Value = 9826
$ws.Range("D3").Value = 96.5

# --- Row 4: adapter name swapped in (was Intel AX200 22.120.1.9, now AX211 23.90.0.2) ---
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 204
$ws.Range("D4").Value = 98.7

# --- Row 5: adapter name swapped in (was Intel AX211 23.90.0.2, now AX200 23.70.2.3) ---
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 25
$ws.Range("D5").Value = 98.90000000000001

# --- Row 6: Totals row ---
$ws.Range("B6").Value = 22
$ws.Range("C6").Value = 10055

# --- Good Drivers table: Total Samples column updates ---
$ws.Range("B16").Value = 338880
$ws.Range("B17").Value = 143869
$ws.Range("B19").Value = 11140
$ws.Range("B21").Value = 14487
$ws.Range("B24").Value = 68450
$ws.Range("B27").Value = 90508
$ws.Range("B30").Value = 52515
